# Registree stats backup on Mon 26 Apr 2021 18:29:12 SAST
#
# This refreshes the "as of" timestamp on every sheet's header cell, inserts
# a newly-registered attendee (Townsend, Diane) into the MD410 Attendance
# sheet's alphabetically-sorted list, and bumps the trailing "Number of
# attendees" / "Number of voters" summary rows to match.

$wb = $excel.ActiveWorkbook

$oldStamp = "24/04/2021 14:19"
$newStamp = "26/04/2021 18:29"

# --- Update the "as of" timestamp header (row 1, col A) on every sheet ---
foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Cells.Item(1, 1)
    $text = $cell.Value2
    if ($text -ne $null -and $text.ToString().Contains($oldStamp)) {
        $cell.Value = $text.ToString().Replace($oldStamp, $newStamp)
    }
}

# --- Insert the new registree row into the "MD410 Attendance" sheet ---
$ws1 = $wb.Worksheets.Item("MD410 Attendance")

# New row goes at row 192 (alphabetically between "Tuckett, Gail" and
# "Toye, Omolayo"), pushing the remaining rows down by one.
$ws1.Rows.Item(192).Insert()

$newRow = $ws1.Range("A192:F192")
$newRow.RowHeight = 25
$newRow.Borders.LineStyle = 1

$ws1.Cells.Item(192, 1).Value = "Townsend"
$ws1.Cells.Item(192, 2).Value = "Diane"
$ws1.Cells.Item(192, 3).Value = "Benoni Lakes"
$ws1.Cells.Item(192, 4).Value = "No"
$ws1.Cells.Item(192, 5).Value = "Yes"
$ws1.Cells.Item(192, 6).Value = "410E"

# --- Refresh the trailing summary rows (now shifted down to 244 / 245) ---
$ws1.Cells.Item(244, 1).Value = "Number of attendees: 241"
$ws1.Cells.Item(245, 1).Value = "Number of voters: 96"
